$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(31810, "Amanda Moraes", "Vendas", "Consulta médica", 2, 45089, 9982.139999999999),
    @(41407, "Matheus da Costa", "Atendimento ao Cliente", "Viagem de negócios", 3, 45087, 7947.9),
    @(47687, "Ana Vitória Souza", "Recursos Humanos", "Problemas pessoais", 1, 45081, 3070.24),
    @(94801, "Ana Vitória Farias", "P&D", "Doença", 6, 45105, 7250.65),
    @(89624, "Nicolas Fogaça", "Marketing", "Problemas pessoais", 2, 45081, 2969.94),
    @(92406, "Bianca Caldeira", "Financeiro", "Viagem de negócios", 4, 45105, 12486.16),
    @(24055, "Lívia Costela", "Vendas", "Problemas pessoais", 8, 45106, 9525.790000000001),
    @(22158, "Carlos Eduardo da Mota", "Vendas", "Outros", 5, 45084, 8287.58),
    @(35994, "Sr. Otávio Cunha", "Atendimento ao Cliente", "Consulta médica", 4, 45094, 5976.63),
    @(45487, "Vitor Hugo Costa", "Recursos Humanos", "Outros", 4, 45082, 4787.74)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $row++
}
